$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (OWASP Top 10 2017 list): translate from Spanish to English ---
$ws.Range("C7").Value  = "A01:2017-Injection"
$ws.Range("C8").Value  = "A02:2017-Broken Authentication"
$ws.Range("C9").Value  = "A03:2017-Sensitive Data Exposure"
$ws.Range("C10").Value = "A04:2017-XML External Entities (XXE)"
$ws.Range("C11").Value = "A05:2017-Broken Access Control"
$ws.Range("C12").Value = "A06:2017-Security Misconfiguration"
$ws.Range("C13").Value = "A07:2017-Cross-Site Scripting (XSS)"
$ws.Range("C14").Value = "A08:2017-Insecure Deserialization"
$ws.Range("C15").Value = "A09:2017-Using Components with Known Vulnerabilities"
$ws.Range("C16").Value = "A10:2017-Insufficient Logging & Monitoring"

# --- Column D "(Nuevo)" badges -> "(New)" ---
$ws.Range("D10").Value = "(New)"
$ws.Range("D14").Value = "(New)"
$ws.Range("D16").Value = "(New)"

# --- Column E (OWASP Top 10 2021 list): corrected / reordered mapping ---
$ws.Range("E7").Value  = "A01:2021-Pérdida de Control de Acceso"
$ws.Range("E8").Value  = "A02:2021-Fallas Criptográficas"
$ws.Range("E9").Value  = "A03:2021-Inyección"
$ws.Range("E10").Value = "A04:2021-Diseño Inseguro"
$ws.Range("E11").Value = "A05:2021-Configuración de Seguridad Incorrecta"
$ws.Range("E12").Value = "A06:2021-Componentes Vulnerables y Desactualizados"
$ws.Range("E13").Value = "A07:2021-Fallas de Identificación y Autenticación"
$ws.Range("E14").Value = "A08:2021-Fallas en la Integridad del Software y de los Datos"
$ws.Range("E15").Value = "A09:2021-Fallas en el Registro y Monitoreo de la Seguridad*"
$ws.Range("E16").Value = "A10:2021-Falsificación de Solicitudes del Lado del Servidor (SSRF)*"

# --- Footnote ---
$ws.Range("E17").Value = "* From the Survey"

# --- Column widths: narrow column C to match the new (shorter) English text ---
$ws.Columns.Item(3).ColumnWidth = 47.64

# --- Remove the two duplicate "New" arrow connectors from the drawing ---
$ws.Shapes.Item("Straight Arrow Connector 12").Delete()
$ws.Shapes.Item("Straight Arrow Connector 11").Delete()

# --- Selection / view state ---
$ws.Range("D17").Select()
